# Insert a new data row at row 200 (pushing existing rows 200-221 down to 201-222),
# matching the "Fruta / hortaliza, semanal" weekly price update for
# "Hortaliza, Vega Monumental Concepción - Ají".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 200..221 down by inserting a fresh row at 200.
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the new weekly record.
$ws.Cells.Item(200, 1).Value = 11
$ws.Cells.Item(200, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(200, 3).Value = "Bíobío"
$ws.Cells.Item(200, 4).Value = 45132
$ws.Cells.Item(200, 5).Value = 8
$ws.Cells.Item(200, 6).Value = 100112021
$ws.Cells.Item(200, 7).Value = "Ají"
$ws.Cells.Item(200, 8).Value = "Inferno"
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 40
$ws.Cells.Item(200, 11).Value = 12000
$ws.Cells.Item(200, 12).Value = 13000
$ws.Cells.Item(200, 13).Value = 12500
$ws.Cells.Item(200, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(200, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(200, 16).Value = 1250
$ws.Cells.Item(200, 17).Value = 10
$ws.Cells.Item(200, 18).Value = "Hortaliza"
